$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E17").Value = 118

$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = 6

$ws.Range("E32").Value = 21
$ws.Range("G32").Value = 9
$ws.Range("H32").Value = 14

$ws.Range("E35").Value = 8

$ws.Range("E52").Value = 6
$ws.Range("F52").Value = 2
$ws.Range("H52").Value = 2

$ws.Range("E56").Value = 7
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 4

$ws.Range("E70").Value = 44
